$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.316374778747559
$ws.Range("B1").Value = 2.3325355052948
$ws.Range("C1").Value = 3.107927322387695
$ws.Range("D1").Value = 3.628358840942383
$ws.Range("E1").Value = 1.919169425964355
